$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2023-05-31 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-01 Thursday", 2) | Out-Null

$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "76+10=86"
$tbl.Cell(1,2).Range.Text = "72-35=37"
$tbl.Cell(1,3).Range.Text = "82-16=66"
$tbl.Cell(1,4).Range.Text = "55+30=85"
$tbl.Cell(1,5).Range.Text = "16+26=42"
$tbl.Cell(2,1).Range.Text = "89-53=36"
$tbl.Cell(2,2).Range.Text = "34-27=7"
$tbl.Cell(2,3).Range.Text = "42+34=76"
$tbl.Cell(2,4).Range.Text = "91+2=93"
$tbl.Cell(2,5).Range.Text = "15+67=82"
$tbl.Cell(3,1).Range.Text = "54-17=37"
$tbl.Cell(3,2).Range.Text = "88+8=96"
$tbl.Cell(3,3).Range.Text = "5-0=5"
$tbl.Cell(3,4).Range.Text = "73-33=40"
$tbl.Cell(3,5).Range.Text = "66+18=84"
$tbl.Cell(4,1).Range.Text = "49-31=18"
$tbl.Cell(4,2).Range.Text = "26+57=83"
$tbl.Cell(4,3).Range.Text = "81-37=44"
$tbl.Cell(4,4).Range.Text = "69+23=92"
$tbl.Cell(4,5).Range.Text = "49-13=36"
$tbl.Cell(5,1).Range.Text = "20+54=74"
$tbl.Cell(5,2).Range.Text = "46-31=15"
$tbl.Cell(5,3).Range.Text = "31+48=79"
$tbl.Cell(5,4).Range.Text = "36+4=40"
$tbl.Cell(5,5).Range.Text = "51+16=67"
$tbl.Cell(6,1).Range.Text = "40+52=92"
$tbl.Cell(6,2).Range.Text = "73-37=36"
$tbl.Cell(6,3).Range.Text = "4+58=62"
$tbl.Cell(6,4).Range.Text = "82-61=21"
$tbl.Cell(6,5).Range.Text = "0+45=45"
$tbl.Cell(7,1).Range.Text = "0+65=65"
$tbl.Cell(7,2).Range.Text = "56+7=63"
$tbl.Cell(7,3).Range.Text = "41+26=67"
$tbl.Cell(7,4).Range.Text = "6+83=89"
$tbl.Cell(7,5).Range.Text = "5+48=53"
$tbl.Cell(8,1).Range.Text = "72+17=89"
$tbl.Cell(8,2).Range.Text = "1+52=53"
$tbl.Cell(8,3).Range.Text = "34+30=64"
$tbl.Cell(8,4).Range.Text = "69-6=63"
$tbl.Cell(8,5).Range.Text = "87-43=44"
$tbl.Cell(9,1).Range.Text = "5+26=31"
$tbl.Cell(9,2).Range.Text = "75-46=29"
$tbl.Cell(9,3).Range.Text = "17+78=95"
$tbl.Cell(9,4).Range.Text = "76-0=76"
$tbl.Cell(9,5).Range.Text = "24+68=92"
$tbl.Cell(10,1).Range.Text = "17+69=86"
$tbl.Cell(10,2).Range.Text = "6+19=25"
$tbl.Cell(10,3).Range.Text = "10+70=80"
$tbl.Cell(10,4).Range.Text = "43-21=22"
$tbl.Cell(10,5).Range.Text = "83-75=8"
$tbl.Cell(11,1).Range.Text = "13+36=49"
$tbl.Cell(11,2).Range.Text = "67-9=58"
$tbl.Cell(11,3).Range.Text = "70+6=76"
$tbl.Cell(11,4).Range.Text = "27-3=24"
$tbl.Cell(11,5).Range.Text = "36+9=45"
$tbl.Cell(12,1).Range.Text = "76-30=46"
$tbl.Cell(12,2).Range.Text = "63+22=85"
$tbl.Cell(12,3).Range.Text = "1+12=13"
$tbl.Cell(12,4).Range.Text = "86-8=78"
$tbl.Cell(12,5).Range.Text = "84-25=59"
$tbl.Cell(13,1).Range.Text = "69-33=36"
$tbl.Cell(13,2).Range.Text = "70+18=88"
$tbl.Cell(13,3).Range.Text = "3+8=11"
$tbl.Cell(13,4).Range.Text = "37+46=83"
$tbl.Cell(13,5).Range.Text = "0+63=63"
$tbl.Cell(14,1).Range.Text = "62-42=20"
$tbl.Cell(14,2).Range.Text = "92-32=60"
$tbl.Cell(14,3).Range.Text = "84-40=44"
$tbl.Cell(14,4).Range.Text = "73+23=96"
$tbl.Cell(14,5).Range.Text = "67-21=46"
$tbl.Cell(15,1).Range.Text = "53-52=1"
$tbl.Cell(15,2).Range.Text = "34-16=18"
$tbl.Cell(15,3).Range.Text = "27+50=77"
$tbl.Cell(15,4).Range.Text = "85-1=84"
$tbl.Cell(15,5).Range.Text = "14-13=1"
$tbl.Cell(16,1).Range.Text = "25+0=25"
$tbl.Cell(16,2).Range.Text = "77-7=70"
$tbl.Cell(16,3).Range.Text = "87-58=29"
$tbl.Cell(16,4).Range.Text = "47-27=20"
$tbl.Cell(16,5).Range.Text = "76-1=75"
$tbl.Cell(17,1).Range.Text = "52-43=9"
$tbl.Cell(17,2).Range.Text = "83+8=91"
$tbl.Cell(17,3).Range.Text = "59-32=27"
$tbl.Cell(17,4).Range.Text = "83-29=54"
$tbl.Cell(17,5).Range.Text = "91-18=73"
$tbl.Cell(18,1).Range.Text = "10+68=78"
$tbl.Cell(18,2).Range.Text = "7+12=19"
$tbl.Cell(18,3).Range.Text = "5+4=9"
$tbl.Cell(18,4).Range.Text = "84-49=35"
$tbl.Cell(18,5).Range.Text = "69-57=12"
$tbl.Cell(19,1).Range.Text = "97+0=97"
$tbl.Cell(19,2).Range.Text = "12+12=24"
$tbl.Cell(19,3).Range.Text = "45-34=11"
$tbl.Cell(19,4).Range.Text = "34+36=70"
$tbl.Cell(19,5).Range.Text = "19-17=2"
$tbl.Cell(20,1).Range.Text = "33+9=42"
$tbl.Cell(20,2).Range.Text = "84-16=68"
$tbl.Cell(20,3).Range.Text = "76-5=71"
$tbl.Cell(20,4).Range.Text = "11-8=3"
$tbl.Cell(20,5).Range.Text = "14+54=68"
